$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 12630.223
$ws.Range("I31").Value = 13584
$ws.Range("K31").Value = 40752
$ws.Range("M31").Value = -40522

$ws.Range("H38").Value = 399.5
$ws.Range("I38").Value = 81.27273
$ws.Range("J38").Value = 3900
$ws.Range("K38").Value = 243.81819
$ws.Range("L38").Value = 11700
$ws.Range("M38").Value = 128.18181
$ws.Range("N38").Value = -12444

$ws.Range("H111").Value = 3835.1667
$ws.Range("I111").Value = 1943
$ws.Range("J111").Value = 5727.3335
$ws.Range("K111").Value = 5829
$ws.Range("L111").Value = 17182.0005
$ws.Range("M111").Value = -2762
$ws.Range("N111").Value = -23316.0005

$ws.Range("H135").Value = 990.1539
$ws.Range("I135").Value = 716
$ws.Range("K135").Value = 6444
$ws.Range("M135").Value = -3909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6002
$ws.Range("I61").Value = 9506
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 9506
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -9294
$ws.Range("N61").Value = -4674

$ws.Range("H63").Value = 4345.909
$ws.Range("I63").Value = 3534.1667
$ws.Range("J63").Value = 5320
$ws.Range("K63").Value = 3534.1667
$ws.Range("L63").Value = 5320
$ws.Range("M63").Value = -2848.1667
$ws.Range("N63").Value = -6692

$ws.Range("H66").Value = 4345.909
$ws.Range("I66").Value = 3534.1667
$ws.Range("J66").Value = 5320
$ws.Range("K66").Value = 17670.8335
$ws.Range("L66").Value = 26600
$ws.Range("M66").Value = -14238.8335
$ws.Range("N66").Value = -33464

$ws.Range("H92").Value = 157725
$ws.Range("J92").Value = 157725
$ws.Range("L92").Value = 157725
$ws.Range("N92").Value = -162717

$ws.Range("H136").Value = 6002
$ws.Range("I136").Value = 9506
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 28518
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -25968
$ws.Range("N136").Value = -17850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 19800
$ws.Range("J19").Value = 19800
$ws.Range("L19").Value = 19800
$ws.Range("N19").Value = -20146

$ws.Range("H35").Value = 39900
$ws.Range("J35").Value = 39900
$ws.Range("L35").Value = 39900
$ws.Range("N35").Value = -40520

$ws.Range("H82").Value = 45143.75
$ws.Range("I82").Value = 44966.668
$ws.Range("J82").Value = 47800
$ws.Range("K82").Value = 44966.668
$ws.Range("L82").Value = 47800
$ws.Range("M82").Value = -44583.668
$ws.Range("N82").Value = -48566

$ws.Range("H85").Value = 45143.75
$ws.Range("I85").Value = 44966.668
$ws.Range("J85").Value = 47800
$ws.Range("K85").Value = 44966.668
$ws.Range("L85").Value = 47800
$ws.Range("M85").Value = -43640.668
$ws.Range("N85").Value = -50452

$ws.Range("H107").Value = 4390.8335
$ws.Range("J107").Value = 3781.6667
$ws.Range("L107").Value = 3781.6667
$ws.Range("N107").Value = -7621.6667

$ws.Range("H112").Value = 25943.32
$ws.Range("J112").Value = 25943.32
$ws.Range("L112").Value = 25943.32
$ws.Range("N112").Value = -28897.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 50000
$ws.Range("I11").Value = 50000
$ws.Range("K11").Value = 50000
$ws.Range("M11").Value = -49860

$ws.Range("H22").Value = 389.3158
$ws.Range("I22").Value = 230.9375
$ws.Range("J22").Value = 1234
$ws.Range("K22").Value = 230.9375
$ws.Range("L22").Value = 1234
$ws.Range("M22").Value = 119.0625
$ws.Range("N22").Value = -1934

$ws.Range("H31").Value = 2696.0645
$ws.Range("I31").Value = 1330.1538
$ws.Range("J31").Value = 9798.799999999999
$ws.Range("K31").Value = 1330.1538
$ws.Range("L31").Value = 9798.799999999999
$ws.Range("M31").Value = -1035.1538
$ws.Range("N31").Value = -10388.8

$ws.Range("H34").Value = 2696.0645
$ws.Range("I34").Value = 1330.1538
$ws.Range("J34").Value = 9798.799999999999
$ws.Range("K34").Value = 1330.1538
$ws.Range("L34").Value = 9798.799999999999
$ws.Range("M34").Value = -1128.1538
$ws.Range("N34").Value = -10202.8

$ws.Range("H58").Value = 3692
$ws.Range("I58").Value = 5537.3335
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 5537.3335
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -5334.3335
$ws.Range("N58").Value = -3406

$ws.Range("H88").Value = 30655
$ws.Range("J88").Value = 30655
$ws.Range("L88").Value = 30655
$ws.Range("N88").Value = -31467

$ws.Range("H91").Value = 30655
$ws.Range("J91").Value = 30655
$ws.Range("L91").Value = 30655
$ws.Range("N91").Value = -33463

$ws.Range("H136").Value = 3692
$ws.Range("I136").Value = 5537.3335
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 16612.0005
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -14062.0005
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6362.3
$ws.Range("I3").Value = 3947.5
$ws.Range("K3").Value = 11842.5
$ws.Range("M3").Value = -11730.5

$ws.Range("H4").Value = 87643.75
$ws.Range("J4").Value = 175
$ws.Range("L4").Value = 525
$ws.Range("N4").Value = -749

$ws.Range("H5").Value = 3834
$ws.Range("I5").Value = 3834
$ws.Range("K5").Value = 11502
$ws.Range("M5").Value = -11390

$ws.Range("H23").Value = 280
$ws.Range("I23").Value = 98
$ws.Range("J23").Value = 362.72726
$ws.Range("K23").Value = 294
$ws.Range("L23").Value = 1088.18178
$ws.Range("M23").Value = -59
$ws.Range("N23").Value = -1558.18178

$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 300
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = -612
$ws.Range("N26").Value = -1476

$ws.Range("H33").Value = 1227.7693
$ws.Range("J33").Value = 1420.1111
$ws.Range("L33").Value = 8520.6666
$ws.Range("N33").Value = -9086.6666

$ws.Range("H63").Value = 352670.66
$ws.Range("J63").Value = 4300
$ws.Range("L63").Value = 12900
$ws.Range("N63").Value = -14398

$ws.Range("H64").Value = 3409
$ws.Range("J64").Value = 4795.8335
$ws.Range("L64").Value = 14387.5005
$ws.Range("N64").Value = -14927.5005

$ws.Range("H66").Value = 352670.66
$ws.Range("J66").Value = 4300
$ws.Range("L66").Value = 38700
$ws.Range("N66").Value = -46188

$ws.Range("H67").Value = 3409
$ws.Range("J67").Value = 4795.8335
$ws.Range("L67").Value = 14387.5005
$ws.Range("N67").Value = -16259.5005

$ws.Range("H118").Value = 3315.4
$ws.Range("I118").Value = 1000
$ws.Range("J118").Value = 3572.6667
$ws.Range("K118").Value = 3000
$ws.Range("L118").Value = 10718.0001
$ws.Range("M118").Value = -1757
$ws.Range("N118").Value = -13204.0001

$ws.Range("H122").Value = 915.1111
$ws.Range("J122").Value = 977.13043
$ws.Range("L122").Value = 8794.173870000001
$ws.Range("N122").Value = -13694.17387

$ws.Range("H131").Value = 18167.623
$ws.Range("I131").Value = 661.2857
$ws.Range("J131").Value = 20831.63
$ws.Range("K131").Value = 1983.8571
$ws.Range("L131").Value = 62494.89
$ws.Range("M131").Value = 3056.1429
$ws.Range("N131").Value = -72574.89

$ws.Range("H132").Value = 1866.64
$ws.Range("I132").Value = 1299.5
$ws.Range("J132").Value = 2133.5293
$ws.Range("K132").Value = 11695.5
$ws.Range("L132").Value = 19201.7637
$ws.Range("M132").Value = -9165.5
$ws.Range("N132").Value = -24261.7637

$ws.Range("H133").Value = 4693.6313
$ws.Range("I133").Value = 1876.2
$ws.Range("J133").Value = 5699.857
$ws.Range("K133").Value = 5628.6
$ws.Range("L133").Value = 17099.571
$ws.Range("M133").Value = -568.6000000000004
$ws.Range("N133").Value = -27219.571

$ws.Range("H135").Value = 3834
$ws.Range("I135").Value = 3834
$ws.Range("K135").Value = 34506
$ws.Range("M135").Value = -31971

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4088.8462
$ws.Range("I132").Value = 3670.6667
$ws.Range("J132").Value = 4143.391
$ws.Range("K132").Value = 11012.0001
$ws.Range("L132").Value = 12430.173
$ws.Range("M132").Value = -8482.000100000001
$ws.Range("N132").Value = -17490.173

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 50000
$ws.Range("K34").Value = 50000
$ws.Range("M34").Value = -49828

$ws.Range("H136").Value = 3154.2
$ws.Range("I136").Value = 2213
$ws.Range("K136").Value = 6639
$ws.Range("M136").Value = -4089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3139.5
$ws.Range("I136").Value = 2546.5454
$ws.Range("J136").Value = 4071.2856
$ws.Range("K136").Value = 7639.6362
$ws.Range("L136").Value = 12213.8568
$ws.Range("M136").Value = -5089.6362
$ws.Range("N136").Value = -17313.8568
